$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035028059928562
$ws.Cells.Item(2, 4).Value = 1.057011448613013
$ws.Cells.Item(2, 5).Value = 1.04566472178648
$ws.Cells.Item(2, 6).Value = 1.060281771490329
$ws.Cells.Item(2, 9).Value = 1.045283661352341
$ws.Cells.Item(2, 10).Value = 1.040144076187173
$ws.Cells.Item(2, 11).Value = 1.059747587309296
$ws.Cells.Item(2, 12).Value = 1.048432325941727
$ws.Cells.Item(2, 13).Value = 1.063008976688834
$ws.Cells.Item(2, 14).Value = 1.017343175891039
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.035833833300696
$ws.Cells.Item(3, 4).Value = 1.057665588844467
$ws.Cells.Item(3, 5).Value = 1.046374937793173
$ws.Cells.Item(3, 6).Value = 1.061060692479221
$ws.Cells.Item(3, 9).Value = 1.045491370423923
$ws.Cells.Item(3, 10).Value = 1.040593953281818
$ws.Cells.Item(3, 11).Value = 1.060215885651222
$ws.Cells.Item(3, 12).Value = 1.04895429187881
$ws.Cells.Item(3, 13).Value = 1.063602382099445
$ws.Cells.Item(3, 14).Value = 1.017493037134781
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.036355889862425
$ws.Cells.Item(4, 4).Value = 1.058089138637627
$ws.Cells.Item(4, 5).Value = 1.046835432628432
$ws.Cells.Item(4, 6).Value = 1.061565586365778
$ws.Cells.Item(4, 9).Value = 1.04562465181476
$ws.Cells.Item(4, 10).Value = 1.040885042955679
$ws.Cells.Item(4, 11).Value = 1.06051849634747
$ws.Cells.Item(4, 12).Value = 1.049292300960236
$ws.Cells.Item(4, 13).Value = 1.063986572322483
$ws.Cells.Item(4, 14).Value = 1.017589977699219
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03657552043244
$ws.Cells.Item(5, 4).Value = 1.058267263319815
$ws.Cells.Item(5, 5).Value = 1.047029247172702
$ws.Cells.Item(5, 6).Value = 1.061778052288148
$ws.Cells.Item(5, 9).Value = 1.045680414191537
$ws.Cells.Item(5, 10).Value = 1.041007413154273
$ws.Cells.Item(5, 11).Value = 1.060645614219436
$ws.Cells.Item(5, 12).Value = 1.049434461252326
$ws.Cells.Item(5, 13).Value = 1.064148136148064
$ws.Cells.Item(5, 14).Value = 1.017630723940533
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03661240658502
$ws.Cells.Item(6, 4).Value = 1.05829717495175
$ws.Cells.Item(6, 5).Value = 1.047061802479321
$ws.Cells.Item(6, 6).Value = 1.061813738418086
$ws.Cells.Item(6, 9).Value = 1.045689761131738
$ws.Cells.Item(6, 10).Value = 1.041027959382734
$ws.Cells.Item(6, 11).Value = 1.060666951978688
$ws.Cells.Item(6, 12).Value = 1.04945833414038
$ws.Cells.Item(6, 13).Value = 1.064175266321325
$ws.Cells.Item(6, 14).Value = 1.017637564956152
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.036358823958261
$ws.Cells.Item(7, 4).Value = 1.058091518498073
$ws.Cells.Item(7, 5).Value = 1.046838021515973
$ws.Cells.Item(7, 6).Value = 1.061568424529114
$ws.Cells.Item(7, 9).Value = 1.045625397972983
$ws.Cells.Item(7, 10).Value = 1.040886678089211
$ws.Cells.Item(7, 11).Value = 1.060520195296317
$ws.Cells.Item(7, 12).Value = 1.049294200274139
$ws.Cells.Item(7, 13).Value = 1.063988730950962
$ws.Cells.Item(7, 14).Value = 1.017590522182829
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035300235580517
$ws.Cells.Item(8, 4).Value = 1.05723245971859
$ws.Cells.Item(8, 5).Value = 1.04590454730874
$ws.Cells.Item(8, 6).Value = 1.060544828016361
$ws.Cells.Item(8, 9).Value = 1.045354089007794
$ws.Cells.Item(8, 10).Value = 1.040296115734756
$ws.Cells.Item(8, 11).Value = 1.059905934787954
$ws.Cells.Item(8, 12).Value = 1.048608671578133
$ws.Cells.Item(8, 13).Value = 1.063209474791101
$ws.Cells.Item(8, 14).Value = 1.017393828065536
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033440056459146
$ws.Cells.Item(9, 4).Value = 1.055720902557661
$ws.Cells.Item(9, 5).Value = 1.044266913243713
$ws.Cells.Item(9, 6).Value = 1.058747953136286
$ws.Cells.Item(9, 9).Value = 1.044867468968083
$ws.Cells.Item(9, 10).Value = 1.039255443193112
$ws.Cells.Item(9, 11).Value = 1.058820462358437
$ws.Cells.Item(9, 12).Value = 1.047402757038516
$ws.Cells.Item(9, 13).Value = 1.061838067135487
$ws.Cells.Item(9, 14).Value = 1.017047021177872
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032203530947369
$ws.Cells.Item(10, 4).Value = 1.05471480828291
$ws.Cells.Item(10, 5).Value = 1.043180159074576
$ws.Cells.Item(10, 6).Value = 1.057554757340729
$ws.Cells.Item(10, 9).Value = 1.044537369677321
$ws.Cells.Item(10, 10).Value = 1.038561721010026
$ws.Cells.Item(10, 11).Value = 1.058094847135318
$ws.Cells.Item(10, 12).Value = 1.046600298424307
$ws.Cells.Item(10, 13).Value = 1.060925070644591
$ws.Cells.Item(10, 14).Value = 1.01681570539882
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.031668976356069
$ws.Cells.Item(11, 4).Value = 1.054279566108471
$ws.Cells.Item(11, 5).Value = 1.042710793481823
$ws.Cells.Item(11, 6).Value = 1.057039236904476
$ws.Cells.Item(11, 9).Value = 1.044393096042548
$ws.Cells.Item(11, 10).Value = 1.038261362047834
$ws.Cells.Item(11, 11).Value = 1.057780201287652
$ws.Cells.Item(11, 12).Value = 1.046253196001099
$ws.Cells.Item(11, 13).Value = 1.060530057666563
$ws.Cells.Item(11, 14).Value = 1.01671552261301
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.031470551342027
$ws.Cells.Item(12, 4).Value = 1.054117960462503
$ws.Cells.Item(12, 5).Value = 1.042536633411119
$ws.Cells.Item(12, 6).Value = 1.056847923294201
$ws.Cells.Item(12, 9).Value = 1.044339306223258
$ws.Cells.Item(12, 10).Value = 1.038149800772922
$ws.Cells.Item(12, 11).Value = 1.057663261869348
$ws.Cells.Item(12, 12).Value = 1.046124323444267
$ws.Cells.Item(12, 13).Value = 1.060383382202323
$ws.Cells.Item(12, 14).Value = 1.01667830750512
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031513108162721
$ws.Cells.Item(13, 4).Value = 1.054152622552414
$ws.Cells.Item(13, 5).Value = 1.042573983014803
$ws.Cells.Item(13, 6).Value = 1.05688895282082
$ws.Cells.Item(13, 9).Value = 1.044350853365293
$ws.Cells.Item(13, 10).Value = 1.038173730778894
$ws.Cells.Item(13, 11).Value = 1.057688348729564
$ws.Cells.Item(13, 12).Value = 1.046151964453083
$ws.Cells.Item(13, 13).Value = 1.060414842298558
$ws.Cells.Item(13, 14).Value = 1.016686290389418
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.031652571772016
$ws.Cells.Item(14, 4).Value = 1.054266206444717
$ws.Cells.Item(14, 5).Value = 1.042696393607083
$ws.Cells.Item(14, 6).Value = 1.057023419307118
$ws.Cells.Item(14, 9).Value = 1.044388653837332
$ws.Cells.Item(14, 10).Value = 1.038252140239953
$ws.Cells.Item(14, 11).Value = 1.057770536374688
$ws.Cells.Item(14, 12).Value = 1.046242542189059
$ws.Cells.Item(14, 13).Value = 1.060517932397764
$ws.Cells.Item(14, 14).Value = 1.016712446454161
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03173851749736
$ws.Cells.Item(15, 4).Value = 1.054336197597162
$ws.Cells.Item(15, 5).Value = 1.042771839149804
$ws.Cells.Item(15, 6).Value = 1.057106291620481
$ws.Cells.Item(15, 9).Value = 1.044411917459159
$ws.Cells.Item(15, 10).Value = 1.038300451653816
$ws.Cells.Item(15, 11).Value = 1.057821166205898
$ws.Cells.Item(15, 12).Value = 1.04629835768576
$ws.Cells.Item(15, 13).Value = 1.060581456281765
$ws.Cells.Item(15, 14).Value = 1.016728561717258
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.032239026038871
$ws.Cells.Item(16, 4).Value = 1.05474370254765
$ws.Cells.Item(16, 5).Value = 1.04321133490588
$ws.Cells.Item(16, 6).Value = 1.0575889949644
$ws.Cells.Item(16, 9).Value = 1.04454691652206
$ws.Cells.Item(16, 10).Value = 1.038581655532809
$ws.Cells.Item(16, 11).Value = 1.058115719808023
$ws.Cells.Item(16, 12).Value = 1.046623342341616
$ws.Cells.Item(16, 13).Value = 1.060951293260827
$ws.Cells.Item(16, 14).Value = 1.016822353787259
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032553215771735
$ws.Cells.Item(17, 4).Value = 1.054999429082986
$ws.Cells.Item(17, 5).Value = 1.043487343312213
$ws.Cells.Item(17, 6).Value = 1.057892089094107
$ws.Cells.Item(17, 9).Value = 1.044631240173203
$ws.Cells.Item(17, 10).Value = 1.038758055538154
$ws.Cells.Item(17, 11).Value = 1.058300366253669
$ws.Cells.Item(17, 12).Value = 1.046827296089925
$ws.Cells.Item(17, 13).Value = 1.061183369304109
$ws.Cells.Item(17, 14).Value = 1.016881181645929
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.032736561064279
$ws.Cells.Item(18, 4).Value = 1.055148628815933
$ws.Cells.Item(18, 5).Value = 1.04364845058636
$ws.Cells.Item(18, 6).Value = 1.058068988736999
$ws.Cells.Item(18, 9).Value = 1.04468029552812
$ws.Cells.Item(18, 10).Value = 1.038860949224925
$ws.Cells.Item(18, 11).Value = 1.058408023906721
$ws.Cells.Item(18, 12).Value = 1.046946294211829
$ws.Cells.Item(18, 13).Value = 1.061318766140186
$ws.Cells.Item(18, 14).Value = 1.016915492831828
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.032799091219869
$ws.Cells.Item(19, 4).Value = 1.055199508581939
$ws.Cells.Item(19, 5).Value = 1.043703403658929
$ws.Cells.Item(19, 6).Value = 1.058129325529383
$ws.Cells.Item(19, 9).Value = 1.044697000193151
$ws.Cells.Item(19, 10).Value = 1.038896033705364
$ws.Cells.Item(19, 11).Value = 1.058444724957616
$ws.Cells.Item(19, 12).Value = 1.046986875445732
$ws.Cells.Item(19, 13).Value = 1.061364938133196
$ws.Cells.Item(19, 14).Value = 1.016927191687439
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032519497512652
$ws.Cells.Item(20, 4).Value = 1.054971988020541
$ws.Cells.Item(20, 5).Value = 1.043457718184604
$ws.Cells.Item(20, 6).Value = 1.057859558572535
$ws.Cells.Item(20, 9).Value = 1.04462220640101
$ws.Cells.Item(20, 10).Value = 1.038739129213982
$ws.Cells.Item(20, 11).Value = 1.058280559924725
$ws.Cells.Item(20, 12).Value = 1.046805410115341
$ws.Cells.Item(20, 13).Value = 1.061158466541654
$ws.Cells.Item(20, 14).Value = 1.016874870188138
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031611499546909
$ws.Cells.Item(21, 4).Value = 1.054232757072507
$ws.Cells.Item(21, 5).Value = 1.042660341665144
$ws.Cells.Item(21, 6).Value = 1.056983817467785
$ws.Cells.Item(21, 9).Value = 1.044377528058671
$ws.Cells.Item(21, 10).Value = 1.038229050447901
$ws.Cells.Item(21, 11).Value = 1.057746335968991
$ws.Cells.Item(21, 12).Value = 1.046215867724865
$ws.Cells.Item(21, 13).Value = 1.06048757353996
$ws.Cells.Item(21, 14).Value = 1.016704744217423
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031041371880689
$ws.Cells.Item(22, 4).Value = 1.053768336894492
$ws.Cells.Item(22, 5).Value = 1.042160060378134
$ws.Cells.Item(22, 6).Value = 1.056434209941952
$ws.Cells.Item(22, 9).Value = 1.044222531422413
$ws.Cells.Item(22, 10).Value = 1.037908375673625
$ws.Cells.Item(22, 11).Value = 1.057410067449947
$ws.Cells.Item(22, 12).Value = 1.045845528050205
$ws.Cells.Item(22, 13).Value = 1.0600660463248
$ws.Cells.Item(22, 14).Value = 1.016597763529362
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03134353407268
$ws.Cells.Item(23, 4).Value = 1.054014499651252
$ws.Cells.Item(23, 5).Value = 1.042425167718217
$ws.Cells.Item(23, 6).Value = 1.056725471238329
$ws.Cells.Item(23, 9).Value = 1.044304807496037
$ws.Cells.Item(23, 10).Value = 1.038078367968425
$ws.Cells.Item(23, 11).Value = 1.057588365302379
$ws.Cells.Item(23, 12).Value = 1.046041820419817
$ws.Cells.Item(23, 13).Value = 1.060289477806875
$ws.Cells.Item(23, 14).Value = 1.016654477350041
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032534733084037
$ws.Cells.Item(24, 4).Value = 1.054984387335948
$ws.Cells.Item(24, 5).Value = 1.043471104145673
$ws.Cells.Item(24, 6).Value = 1.057874257375486
$ws.Cells.Item(24, 9).Value = 1.044626288773141
$ws.Cells.Item(24, 10).Value = 1.038747681197496
$ws.Cells.Item(24, 11).Value = 1.058289509686818
$ws.Cells.Item(24, 12).Value = 1.046815299336154
$ws.Cells.Item(24, 13).Value = 1.061169718933125
$ws.Cells.Item(24, 14).Value = 1.016877722070817
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.033920331610462
$ws.Cells.Item(25, 4).Value = 1.056111401359218
$ws.Cells.Item(25, 5).Value = 1.044689407779608
$ws.Cells.Item(25, 6).Value = 1.059211665557939
$ws.Cells.Item(25, 9).Value = 1.044994278139955
$ws.Cells.Item(25, 10).Value = 1.039524477145652
$ws.Cells.Item(25, 11).Value = 1.059101436641502
$ws.Cells.Item(25, 12).Value = 1.04771425976605
$ws.Cells.Item(25, 13).Value = 1.062192391757575
$ws.Cells.Item(25, 14).Value = 1.017136700702212

Write-Host "Updated 264 cells in vm_pu data range (rows 2-25)"
